$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price values in column B (base price)
$ws.Range("B9").Value = 5500
$ws.Range("B10").Value = 5300
$ws.Range("B11").Value = 5100
$ws.Range("B12").Value = 5000
$ws.Range("B13").Value = 4900

# Update price per km values in column D
$ws.Range("D8").Value = 25
$ws.Range("D9").Value = 24
$ws.Range("D10").Value = 23.75
$ws.Range("D11").Value = 23.5
$ws.Range("D12").Value = 23.25
$ws.Range("D13").Value = 23

# Update view: scroll back to top-left A1 and change selection to G11
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("G11").Select()
